# Applies the Fri Aug 25 11:47:33 UTC 2023 "cryptos list" refresh:
# updated Price/Volume(1h) figures for the existing rows, and a swap of
# the #49/#50 ranking rows (Cronos <-> Mantle, incl. link + price + vol).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 cell -> new text value, in the same order as the source diff.
$updates = [ordered]@{
    "D2" = "26.196.96";
    "E2" = "  -1.34%  ";
    "D3" = "1.659.36";
    "E3" = "  -0.85%  ";
    "D4" = "1.004";
    "E4" = "  +0.23%  ";
    "D5" = "217.14";
    "E5" = "  -1.39%  ";
    "D6" = "0.5190";
    "E6" = "  -1.85%  ";
    "E7" = "  +0.20%  ";
    "E8" = "  -1.63%  ";
    "D9" = "0.06277";
    "E9" = "  -1.84%  ";
    "D10" = "20.79";
    "E10" = "  -5.20%  ";
    "D11" = "0.07773";
    "E11" = "  -0.38%  ";
    "D12" = "4.481";
    "E12" = "  -0.23%  ";
    "D13" = "1.650.22";
    "E13" = "  -1.47%  ";
    "D14" = "1.886.53";
    "E14" = "  -0.86%  ";
    "D15" = "0.5470";
    "E15" = "  -2.01%  ";
    "D16" = "0.0₅8121";
    "E16" = "  -2.63%  ";
    "D17" = "65.01";
    "E17" = "  -1.12%  ";
    "D18" = "26.216.17";
    "E18" = "  -1.24%  ";
    "E19" = "  +0.22%  ";
    "E20" = "  -3.20%  ";
    "D21" = "191.99";
    "E21" = "  -0.79%  ";
    "E22" = "  -2.75%  ";
    "D23" = "6.008";
    "E23" = "  -5.08%  ";
    "D25" = "139.40";
    "E25" = "  -0.23%  ";
    "D26" = "0.1224";
    "E26" = "  -3.61%  ";
    "D27" = "7.299";
    "E27" = "  -1.61%  ";
    "E28" = "  -0.96%  ";
    "D29" = "1.442";
    "E29" = "  +1.04%  ";
    "D30" = "0.05936";
    "E30" = "  -4.27%  ";
    "D31" = "1.274";
    "E31" = "  -1.47%  ";
    "D32" = "3.544";
    "E32" = "  -2.09%  ";
    "D33" = "3.281";
    "E33" = "  -4.39%  ";
    "E34" = "  -6.06%  ";
    "D35" = "0.9612";
    "E35" = "  -4.79%  ";
    "E36" = "  +0.11%  ";
    "D37" = "2.770";
    "E37" = "  -0.36%  ";
    "D38" = "0.5697";
    "E38" = "  -6.69%  ";
    "D39" = "6.034";
    "E39" = "  -0.65%  ";
    "D40" = "0.01591";
    "E40" = "  -1.75%  ";
    "D41" = "0.8526";
    "E41" = "  -0.43%  ";
    "E42" = "  +0.24%  ";
    "D43" = "1.010.08";
    "E43" = "  -7.68%  ";
    "D44" = "100.63";
    "E44" = "  -0.14%  ";
    "D45" = "1.801.00";
    "E45" = "  -0.95%  ";
    "D46" = "0.0₈109";
    "E46" = "  -2.67%  ";
    "D47" = "56.52";
    "E47" = "  -3.25%  ";
    "D48" = "1.007";
    "E48" = "  -0.14%  ";
    "D49" = "8.022";
    "E49" = "  -1.27%  ";
    "B50" = "Mantle";
    "C50" = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt";
    "D50" = "0.4269";
    "E50" = "  +0.94%  ";
    "B51" = "Cronos";
    "C51" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro";
    "D51" = "0.05167";
    "E51" = "  -0.72%  "
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # The sheet stores every Price/Volume(1h) entry as literal text (even
    # the numeric-looking ones, e.g. "1.004", "217.14"). Force text format
    # before assigning so Excel doesn't reinterpret/round them as numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
